$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Unhide rows 6-28 (previously hidden diagnostic/module rows)
$ws.Rows("6:28").Hidden = $false

# 2. Insert two new rows after row 110 (pushes everything below down by 2)
$ws.Rows("111:112").Insert()

# 3. Populate the new header row 110 with the new section title
$ws.Range("A110").Value = "End of Project Module Direct Dependencies"
$ws.Range("A110").Font.Bold = $true
$ws.Range("A110").HorizontalAlignment = -4131
$ws.Range("A110").VerticalAlignment = -4160
$ws.Rows("110").RowHeight = 13.9

# 4. Bold the "Adafruit-Blinka" module name label (now row 113, was row 111)
$ws.Range("A113").Font.Bold = $true

# 5. Update the active selection to match where the user left off editing
$ws.Range("E121").Select()

Write-Host "edits applied"
